# Update the nationality for the row 5 student (李立化 / REFAEL CORNELIUS
# HARIANJA) from 印尼 (Indonesia) to 越南 (Vietnam), and leave the active
# selection where the author last left it (D14) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "越南"

$ws.Range("D14").Select()
